# Apply the "cryptos list" data refresh captured by the commit diff.
# Numeric-looking text values (e.g. "1.00", "7.61") are written with a
# leading single-quote so Excel keeps them as literal text instead of
# silently re-parsing them as numbers (which would drop trailing zeros
# or introduce floating point noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.916.17'
$ws.Range("E2").Value = '  +0.96%  '

# Row 3
$ws.Range("D3").Value = '3.375.75'
$ws.Range("E3").Value = '  +0.04%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = '''568.94'
$ws.Range("E5").Value = '  +0.33%  '

# Row 6
$ws.Range("D6").Value = '''139.87'
$ws.Range("E6").Value = '  -0.32%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").Value = '''7.61'
$ws.Range("E9").Value = '  +1.68%  '

# Row 10
$ws.Range("E10").Value = '  -1.40%  '

# Row 11
$ws.Range("D11").Value = '''0.385'
$ws.Range("E11").Value = '  -0.57%  '

# Row 12
$ws.Range("D12").Value = '3.953.06'
$ws.Range("E12").Value = '  +0.05%  '

# Row 13
$ws.Range("E13").Value = '  +1.83%  '

# Row 14
$ws.Range("D14").Value = '''27.73'
$ws.Range("E14").Value = '  -1.23%  '

# Row 15
$ws.Range("D15").Value = '3.381.98'
$ws.Range("E15").Value = '  +0.31%  '

# Row 16
$ws.Range("E16").Value = '  +0.03%  '

# Row 17
$ws.Range("D17").Value = '61.036.60'
$ws.Range("E17").Value = '  +0.91%  '

# Row 18
$ws.Range("D18").Value = '''6.07'
$ws.Range("E18").Value = '  -1.90%  '

# Row 19
$ws.Range("D19").Value = '''13.52'
$ws.Range("E19").Value = '  -2.11%  '

# Row 20
$ws.Range("D20").Value = '''8.86'
$ws.Range("E20").Value = '  -1.46%  '

# Row 21
$ws.Range("D21").Value = '''381.39'
$ws.Range("E21").Value = '  -1.20%  '

# Row 22
$ws.Range("D22").Value = '''75.38'
$ws.Range("E22").Value = '  +3.11%  '

# Row 23
$ws.Range("D23").Value = '''0.548'

# Row 24
$ws.Range("E24").Value = '  -0.15%  '

# Row 25
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.517.38'
$ws.Range("E25").Value = '  -0.05%  '

# Row 26
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = '''0.0000113'
$ws.Range("E26").Value = '  -1.92%  '

# Row 27
$ws.Range("E27").Value = '  +7.07%  '

# Row 28
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  +0.14%  '

# Row 29
$ws.Range("D29").Value = '''7.17'
$ws.Range("E29").Value = '  -2.38%  '

# Row 30
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("E31").Value = '  -0.46%  '

# Row 32
$ws.Range("E32").Value = '  -0.02%  '

# Row 33
$ws.Range("E33").Value = '  -3.59%  '

# Row 34
$ws.Range("D34").Value = '''23.14'
$ws.Range("E34").Value = '  -1.58%  '

# Row 35
$ws.Range("E35").Value = '  +0.00%  '

# Row 36
$ws.Range("D36").Value = '''166.21'
$ws.Range("E36").Value = '  -1.35%  '

# Row 37
$ws.Range("D37").Value = '3.411.26'
$ws.Range("E37").Value = '  +0.21%  '

# Row 38
$ws.Range("D38").Value = '''4.94'
$ws.Range("E38").Value = '  +0.41%  '

# Row 39
$ws.Range("E39").Value = '  -2.87%  '

# Row 40
$ws.Range("D40").Value = '''0.0761'
$ws.Range("E40").Value = '  -1.06%  '

# Row 41
$ws.Range("D41").Value = '''25.88'
$ws.Range("E41").Value = '  -4.76%  '

# Row 43
$ws.Range("E43").Value = '  +0.18%  '

# Row 44
$ws.Range("E44").Value = '  -1.87%  '

# Row 45
$ws.Range("E45").Value = '  -3.46%  '

# Row 46
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("D47").Value = '2.427.18'
$ws.Range("E47").Value = '  -3.55%  '

# Row 48
$ws.Range("D48").Value = '''22.93'
$ws.Range("E48").Value = '  -0.73%  '

# Row 49
$ws.Range("E49").Value = '  -1.99%  '

# Row 50
$ws.Range("E50").Value = '  -3.11%  '

# Row 51
$ws.Range("E51").Value = '  +6.32%  '
